$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of homework grade data
$ws.Range("A3").Value = "hw1_part2"
$ws.Range("B3").Value = "script looks OK"
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = "hw2"
$ws.Range("B4").Value = "script looks OK, code in ulib.s OK"
$ws.Range("C4").Value = 10

# Copy style from A2:B2 down to A3:B4 (already inherited from header copy, but ensure formatting)
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B4").PasteSpecial(-4122)  # xlPasteFormats

# Center align the Scores column header and values
$ws.Range("C1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C2:C4").HorizontalAlignment = -4108  # xlCenter

# Update selection to match new used range
$ws.Range("A1:C4").Select()
